$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date in column C for rows 2-15
# from 2023-09-15 (45184) to 2023-09-16 (45185)
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 3).Value = 45185
}
